# Rename source_data sheet ("Data table") to "Data"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data table")
$ws.Name = "Data"

# Update the selection on that sheet to H26 (single cell)
$ws.Activate()
$ws.Range("H26").Select()
